# add BY2023 election and gov and HE2023 election
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: ABG (Aktion Buerger fuer Gerechtigkeit) - Hessen 2023
$ws.Range("A20").Value = "ABG"
$ws.Range("B20").Value = "Aktion Bürger für Gerechtigkeit"
$ws.Range("H20").Value = "https://www.bpb.de/themen/parteien/wer-steht-zur-wahl/hessen-2023/539442/aktion-buerger-fuer-gerechtigkeit/"

# Row 21: DNM (Die Neue Mitte) - Hessen 2023
$ws.Range("A21").Value = "DNM"
$ws.Range("B21").Value = "Die Neue Mitte"
$ws.Range("H21").Value = "https://www.bpb.de/themen/parteien/wer-steht-zur-wahl/hessen-2023/539446/die-neue-mitte/"

# Row 22: WKH (Waehlerliste Klimaliste Hessen) - Hessen 2023
$ws.Range("A22").Value = "WKH"
$ws.Range("H22").Value = "https://www.bpb.de/themen/parteien/wer-steht-zur-wahl/hessen-2023/539448/waehlerliste-klimaliste-hessen/"
$ws.Range("B22").Value = "Wählerliste Klimaliste Hessen "
$ws.Range("G22").Value = "Not formally affiliated with the Klimaliste"

$ws.Range("B19").Select()
